{"js": "// Learning diary \u2013 bug fix: the project date range end-date was wrong.\n// \"24.6.2024 \u2013 26.6.2024\"  ->  \"24.6.2024 \u2013 27.6.2024\"\n//\n// The original run of text is split (by the author's live edit) into\n// \"24.6.2024 \u2013 2\" + \"7\" + \".6.2024\", but the net, visible effect is simply\n// changing the end day from the 26th to the 27th. We find the unique\n// paragraph containing the old date range and rewrite its text.\n\nconst oldText = \"24.6.2024 \\u2013 26.6.2024\"; // en dash (U+2013)\nconst newText = \"24.6.2024 \\u2013 27.6.2024\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the date range text '\" + oldText + \"' to update.\");\n}\n\n// Replace in place, preserving the surrounding paragraph/formatting.\nresults.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Learning diary - bug fix: the project date range end-date was wrong.\n# \"24.6.2024 - 26.6.2024\"  ->  \"24.6.2024 - 27.6.2024\"\n#\n# The original run of text was split (by the author's live edit) into\n# \"24.6.2024 - 2\" + \"7\" + \".6.2024\", but the net, visible effect is simply\n# changing the end day from the 26th to the 27th. We locate the unique\n# paragraph containing the old date range and replace its text via Find/Replace.\n\n$d = $word.ActiveDocument\n\n$oldText = \"24.6.2024 \" + [char]0x2013 + \" 26.6.2024\"   # en dash (U+2013)\n$newText = \"24.6.2024 \" + [char]0x2013 + \" 27.6.2024\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\nif (-not $found) {\n    throw \"Could not find the date range text '$oldText' to update.\"\n}\n"}
